$d = $word.ActiveDocument

function Replace-Exact($old, $new) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND: $old"
    }
}

# 1. Title
Replace-Exact "Relatório de tendência de mercado: Contoso Protein Plus" "Relatório de tendências de mercado: Contoso Protein Plus"

# 2-3. "Conteúdo envolvente:" section
Replace-Exact "Conteúdo envolvente:" "Conteúdo interessante:"
Replace-Exact " O rolo viral transmitiu com sucesso o apelo do produto por meio de visuais envolventes e conteúdo informativo." " O reel viral transmitiu com sucesso o apelo do produto por meio de recursos visuais interessantes e conteúdo informativo."

# 4-5. "Marketing de Influência:" section
Replace-Exact "Marketing de Influência:" "Marketing de influenciadores:"
Replace-Exact " O poder do marketing de influência não pode ser exagerado." " não dá para superestimar o poder do marketing de influenciadores."

# 6-7. "Variedades de sabor e sabor:" section
Replace-Exact "Variedades de sabor e sabor:" "Variedades de sabores:"
Replace-Exact " a reputação da Contoso Protein Plus por sabores deliciosos e diversos foi um ponto chave de venda no conteúdo viral." " a reputação do Contoso Protein Plus pelos sabores deliciosos e variados foi um ponto de venda importante no conteúdo viral."

# 8-9. "Tendências de saúde e fitness:" section
Replace-Exact "Tendências de saúde e fitness:" "Tendências de saúde e bem-estar:"
Replace-Exact " O aumento contínuo na consciência de saúde e fitness, combinado com um aumento no número de pessoas que adotam rotinas de treino e estilos de vida ativos, criou um mercado receptivo para um produto como o Contoso Protein Plus." " a onda atual da conscientização sobre saúde e bem-estar, combinada com um aumento no número de pessoas aderindo a rotinas de exercícios e estilos de vida ativos, criaram um mercado receptivo para um produto como o Contoso Protein Plus."

# 10-11. Bold "Fácil disponibilidade:" heading (including the space run before it)
$rng = $d.Content
$ok = $rng.Find.Execute(" Fácil disponibilidade:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($ok) {
    $rng.Font.Bold = 1
    $rng.Font.BoldBi = 1
}

# 12. "hype" sentence
Replace-Exact " A acessibilidade do produto através de vários varejistas on-line alimentou ainda mais o hype." " a acessibilidade do produto por meio de vários varejistas on-line estimulou ainda mais o hype."

# 13-14. "Comentários positivos e depoimentos:" section
Replace-Exact "Comentários positivos e depoimentos:" "Críticas e avaliações positivas:"
Replace-Exact " O rolo não foi um caso isolado." " o reel não foi um caso isolado."

# 15. "Boca a boca:" section sentence
Replace-Exact " As plataformas de mídia social promovem a rápida disseminação de tendências através do boca a boca." " as plataformas de rede social promovem a rápida disseminação das tendências por meio do boca a boca."
